# Weekly update: a new price record (week of 2023-05-04) is inserted as
# row 9, pushing the existing rows 9-82 down to rows 10-83 (dimension
# grows from A1:R82 to A1:R83). The new row reuses the same constant
# market/category metadata (columns A, B, C, E, F, G, R) as every other
# data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 9, shifting rows 9:82 down to 10:83.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data point.
$ws.Cells.Item(9, 1).Value  = 1
$ws.Cells.Item(9, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value  = 45050
$ws.Cells.Item(9, 5).Value  = 15
$ws.Cells.Item(9, 6).Value  = 100112031
$ws.Cells.Item(9, 7).Value  = "Poroto verde"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Segunda"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 28000
$ws.Cells.Item(9, 12).Value = 30000
$ws.Cells.Item(9, 13).Value = 29000
$ws.Cells.Item(9, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 1160
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
